$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45171 -> 45172) for every data row (rows 2 through 295).
$ws.Range("C2:C295").Value = 45172
